# Scheduled runner update: refresh market-price-derived profit columns
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the per-job leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 808.4
$ws.Range("I12").Value = 869.4
$ws.Range("J12").Value = 747.4
$ws.Range("K12").Value = 869.4
$ws.Range("L12").Value = 747.4
$ws.Range("M12").Value = -699.4
$ws.Range("N12").Value = -1087.4

# Row 17
$ws.Range("H17").Value = 2458.6365
$ws.Range("J17").Value = 2458.6365
$ws.Range("L17").Value = 7375.9095
$ws.Range("N17").Value = -7711.9095

# Row 40
$ws.Range("H40").Value = 25411.172
$ws.Range("I40").Value = 16119.934
$ws.Range("J40").Value = 35366.07
$ws.Range("K40").Value = 16119.934
$ws.Range("L40").Value = 35366.07
$ws.Range("M40").Value = -15944.934
$ws.Range("N40").Value = -35716.07

# Row 99
$ws.Range("H99").Value = 6799.6
$ws.Range("J99").Value = 13140.429
$ws.Range("L99").Value = 39421.287
$ws.Range("N99").Value = -42417.287

# Row 103
$ws.Range("H103").Value = 519.7
$ws.Range("I103").Value = 489.1111
$ws.Range("J103").Value = 795
$ws.Range("K103").Value = 1467.3333
$ws.Range("L103").Value = 2385
$ws.Range("M103").Value = -881.3333
$ws.Range("N103").Value = -3557

# Row 132
$ws.Range("H132").Value = 9495.441999999999
$ws.Range("I132").Value = 1668.4186
$ws.Range("K132").Value = 5005.2558
$ws.Range("M132").Value = -2475.2558

# Row 135
$ws.Range("H135").Value = 2464.3225
$ws.Range("I135").Value = 850.36
$ws.Range("K135").Value = 7653.24
$ws.Range("M135").Value = -5118.24

# Row 137
$ws.Range("H137").Value = 12719552
$ws.Range("J137").Value = 22227354
$ws.Range("L137").Value = 66682062
$ws.Range("N137").Value = -66687162

# Row 138
$ws.Range("H138").Value = 3627.449
$ws.Range("I138").Value = 1211.3928
$ws.Range("J138").Value = 4593.8716
$ws.Range("K138").Value = 3634.1784
$ws.Range("L138").Value = 13781.6148
$ws.Range("M138").Value = 1505.8216
$ws.Range("N138").Value = -24061.6148

# Row 141
$ws.Range("H141").Value = 2504.3777
$ws.Range("I141").Value = 2480.2856
$ws.Range("J141").Value = 2841.6667
$ws.Range("K141").Value = 7440.8568
$ws.Range("L141").Value = 8525.000100000001
$ws.Range("M141").Value = -2260.8568
$ws.Range("N141").Value = -18885.0001

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 280.33334
$ws.Range("I5").Value = 280.33334
$ws.Range("K5").Value = 280.33334
$ws.Range("M5").Value = -168.33334

# Row 61
$ws.Range("H61").Value = 21678.445
$ws.Range("I61").Value = 33621.2
$ws.Range("J61").Value = 6750
$ws.Range("K61").Value = 33621.2
$ws.Range("L61").Value = 6750
$ws.Range("M61").Value = -33409.2
$ws.Range("N61").Value = -7174

# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 110
$ws.Range("H110").Value = 929809.75
$ws.Range("I110").Value = 1135084.1
$ws.Range("K110").Value = 1135084.1
$ws.Range("M110").Value = -1133039.1

# Row 136
$ws.Range("H136").Value = 21678.445
$ws.Range("I136").Value = 33621.2
$ws.Range("J136").Value = 6750
$ws.Range("K136").Value = 100863.6
$ws.Range("L136").Value = 20250
$ws.Range("M136").Value = -98313.59999999999
$ws.Range("N136").Value = -25350

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 280.33334
$ws.Range("I4").Value = 280.33334
$ws.Range("K4").Value = 280.33334
$ws.Range("M4").Value = -165.33334

# Row 22
$ws.Range("H22").Value = 414.0625
$ws.Range("I22").Value = 482.8889
$ws.Range("J22").Value = 325.57144
$ws.Range("K22").Value = 482.8889
$ws.Range("L22").Value = 325.57144
$ws.Range("M22").Value = -309.8889
$ws.Range("N22").Value = -671.5714399999999

# Row 86
$ws.Range("H86").Value = 4087.4443
$ws.Range("I86").Value = 3430.4
$ws.Range("K86").Value = 3430.4
$ws.Range("M86").Value = -2307.4

# Row 89
$ws.Range("H89").Value = 4087.4443
$ws.Range("I89").Value = 3430.4
$ws.Range("K89").Value = 17152
$ws.Range("M89").Value = -11536

# Row 134
$ws.Range("H134").Value = 1138.1063
$ws.Range("I134").Value = 1101.4147
$ws.Range("J134").Value = 1388.8334
$ws.Range("K134").Value = 3304.2441
$ws.Range("L134").Value = 4166.5002
$ws.Range("M134").Value = -769.2440999999999
$ws.Range("N134").Value = -9236.5002

$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

# Row 31
$ws.Range("H31").Value = 6504.1904
$ws.Range("I31").Value = 1811.75
$ws.Range("K31").Value = 1811.75
$ws.Range("M31").Value = -1516.75

# Row 34
$ws.Range("H34").Value = 6504.1904
$ws.Range("I34").Value = 1811.75
$ws.Range("K34").Value = 1811.75
$ws.Range("M34").Value = -1609.75

# Row 58
$ws.Range("H58").Value = 337391.94
$ws.Range("I58").Value = 1905.2727
$ws.Range("J58").Value = 531621.0600000001
$ws.Range("K58").Value = 1905.2727
$ws.Range("L58").Value = 531621.0600000001
$ws.Range("M58").Value = -1702.2727
$ws.Range("N58").Value = -532027.0600000001

# Row 74
$ws.Range("H74").Value = 58657
$ws.Range("J74").Value = 58657
$ws.Range("L74").Value = 58657
$ws.Range("N74").Value = -60405

# Row 77
$ws.Range("H77").Value = 58657
$ws.Range("J77").Value = 58657
$ws.Range("L77").Value = 175971
$ws.Range("N77").Value = -184707

# Row 134
$ws.Range("H134").Value = 1499.081
$ws.Range("I134").Value = 1580.4546
$ws.Range("K134").Value = 4741.3638
$ws.Range("M134").Value = -2206.3638

# Row 136
$ws.Range("H136").Value = 337391.94
$ws.Range("I136").Value = 1905.2727
$ws.Range("J136").Value = 531621.0600000001
$ws.Range("K136").Value = 5715.8181
$ws.Range("L136").Value = 1594863.18
$ws.Range("M136").Value = -3165.8181
$ws.Range("N136").Value = -1599963.18

# Row 141
$ws.Range("H141").Value = 86238.125
$ws.Range("J141").Value = 87967.60000000001
$ws.Range("L141").Value = 87967.60000000001
$ws.Range("N141").Value = -98327.60000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 2873464
$ws.Range("J11").Value = 6670299.5
$ws.Range("L11").Value = 20010898.5
$ws.Range("N11").Value = -20011178.5

# Row 34
$ws.Range("H34").Value = 2039
$ws.Range("J34").Value = 2000.75
$ws.Range("L34").Value = 6002.25
$ws.Range("N34").Value = -6170.25

# Row 86
$ws.Range("H86").Value = 499.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 499.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 1498.5
$ws.Range("N86").Value = -3870.5
$ws.Range("M86").ClearContents()

# Row 89
$ws.Range("H89").Value = 499.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 499.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 4495.5
$ws.Range("N89").Value = -16351.5
$ws.Range("M89").ClearContents()

# Row 118
$ws.Range("H118").Value = 2846.375
$ws.Range("J118").Value = 4298.6665
$ws.Range("L118").Value = 12895.9995
$ws.Range("N118").Value = -15381.9995

# Row 132
$ws.Range("H132").Value = 3878.8823
$ws.Range("I132").Value = 667.625
$ws.Range("K132").Value = 6008.625
$ws.Range("M132").Value = -3478.625

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 422.7742
$ws.Range("I2").Value = 255.9375
$ws.Range("J2").Value = 600.73334
$ws.Range("K2").Value = 255.9375
$ws.Range("L2").Value = 600.73334
$ws.Range("M2").Value = -142.9375
$ws.Range("N2").Value = -826.73334

# Row 132
$ws.Range("H132").Value = 863212.3
$ws.Range("I132").Value = 289996.72
$ws.Range("K132").Value = 869990.1599999999
$ws.Range("M132").Value = -867460.1599999999

# Row 140
$ws.Range("H140").Value = 63259.332
$ws.Range("J140").Value = 63259.332
$ws.Range("L140").Value = 63259.332
$ws.Range("N140").Value = -73619.33199999999

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3958.611
$ws.Range("I16").Value = 2159.3333
$ws.Range("J16").Value = 7557.1665
$ws.Range("K16").Value = 2159.3333
$ws.Range("L16").Value = 7557.1665
$ws.Range("M16").Value = -1989.3333
$ws.Range("N16").Value = -7897.1665

# Row 132
$ws.Range("H132").Value = 6515.385
$ws.Range("I132").Value = 7495.25
$ws.Range("K132").Value = 22485.75
$ws.Range("M132").Value = -19955.75

# Row 136
$ws.Range("H136").Value = 3387.0417
$ws.Range("I136").Value = 2210.5557
$ws.Range("K136").Value = 6631.6671
$ws.Range("M136").Value = -4081.6671

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 79286
$ws.Range("J46").Value = 79286
$ws.Range("L46").Value = 79286
$ws.Range("N46").Value = -79748

# Row 126
$ws.Range("H126").Value = 4950.8
$ws.Range("I126").Value = 4438.5
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 13315.5
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -10845.5
$ws.Range("N126").Value = -25940

# Row 132
$ws.Range("H132").Value = 17858466
$ws.Range("I132").Value = 1580.6875
$ws.Range("J132").Value = 41667650
$ws.Range("K132").Value = 4742.0625
$ws.Range("L132").Value = 125002950
$ws.Range("M132").Value = -2212.0625
$ws.Range("N132").Value = -125008010

# Row 134
$ws.Range("H134").Value = 79286
$ws.Range("J134").Value = 79286
$ws.Range("L134").Value = 237858
$ws.Range("N134").Value = -242928

# Row 136
$ws.Range("H136").Value = 7237.19
$ws.Range("I136").Value = 2474.2273
$ws.Range("J136").Value = 10979.518
$ws.Range("K136").Value = 7422.6819
$ws.Range("L136").Value = 32938.554
$ws.Range("M136").Value = -4872.6819
$ws.Range("N136").Value = -38038.554
